$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New device-master rows (device set "32"): Finger Print Scanner, IRIS
# Scanner, Web Camera, Document Scanner, Printer - following the existing
# id/mac/serial pattern used for sets 1..31.
$names   = @("Finger Print Scanner 32", "IRIS Scanner 32", "Web Camera 32", "Document Scanner 32", "Printer 32")
$macs    = @("80-75-40-E8-CA-24", "0E-1A-14-4A-6D-3A", "65-13-7F-0F-F7-53", "73-C4-DE-8E-C9-8D", "EC-74-AB-E0-0F-38")
$serials = @("BS563Q2230824", "BS563Q2230825", "BS563Q2230826", "BS563Q2230827", "BS563Q2230828")
$dspecs  = @(165, 327, 736, 801, 920)

$startId  = 3000176
$startRow = 157

for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $startId + $i
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $names[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 3).Value = $macs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $serials[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $dspecs[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $ws.Cells.Item($startRow + $i, 7).Value = "eng"
    $ws.Cells.Item($startRow + $i, 8).Value = $true
    $ws.Cells.Item($startRow + $i, 8).HorizontalAlignment = -4131
    $ws.Cells.Item($startRow + $i, 9).Value = "superadmin"
    $ws.Cells.Item($startRow + $i, 10).Value = "now()"
}

# Scroll the view down (topLeftCell A113) and select from K113 to the end
# of the grid (K1:XFD1048576 in the authored edit), as in the original
# session.
$win = $excel.ActiveWindow
$win.ScrollRow = 113
$win.ScrollColumn = 1
$ws.Range("K113:XFD1048576").Select()
